$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Range("A8").Value  = "{'Student': np.int64(1), 'Hobby': np.int64(1)}"
$ws.Range("A12").Value = "{'Student': np.int64(1), 'Gender': np.int64(1)}"
$ws.Range("A13").Value = "{'Student': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A16").Value = "{'HDI': np.int64(1), 'Gender': np.int64(1)}"
$ws.Range("A17").Value = "{'HDI': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A18").Value = "{'Gender': np.int64(1), 'SexualOrientation': np.int64(1), 'Hobby': np.int64(1)}"
$ws.Range("A19").Value = "{'Student': np.int64(1), 'SexualOrientation': np.int64(1), 'Gender': np.int64(1)}"
$ws.Range("A20").Value = "{'HDI': np.int64(1), 'Student': np.int64(1), 'Gender': np.int64(1)}"
$ws.Range("A21").Value = "{'Student': np.int64(1), 'HDI': np.int64(1), 'SexualOrientation': np.int64(1)}"
$ws.Range("A22").Value = "{'HDI': np.int64(1), 'SexualOrientation': np.int64(1), 'Gender': np.int64(1)}"
